# Apply the commit's changes:
#  1. Rename sheet "Пункт 1" -> "Графики"
#  2. Add new defined names Y_linear / Y_degree / Y_poly / Y_exp
#  3. On "Графики" sheet, add a "corr" column (E) with CORREL() formulas
#     and an F column that flags the best-fitting approximation
#  4. Re-point the active tab at "Графики" and fix up sheet selections

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Пункт 1" sheet to "Графики" ---------------------------
$wsGraphs = $wb.Worksheets.Item("Пункт 1")
$wsGraphs.Name = "Графики"

# --- 2. New named ranges on "Данные" for each approximation's Y column ----
$wb.Names.Add("Y_linear", "=Данные!`$D`$2:`$D`$21")
$wb.Names.Add("Y_degree", "=Данные!`$E`$2:`$E`$21")
$wb.Names.Add("Y_poly",   "=Данные!`$F`$2:`$F`$21")
$wb.Names.Add("Y_exp",    "=Данные!`$G`$2:`$G`$21")

# --- 3. Correlation column on "Графики" ------------------------------------
$wsGraphs.Range("E1").Value = "corr"
$wsGraphs.Range("E1").Font.Italic = $true
$wsGraphs.Range("E1").Font.Family = 1

$wsGraphs.Range("E2").Formula = "=CORREL(Y,Y_linear)"
$wsGraphs.Range("E3").Formula = "=CORREL(Y,Y_degree)"
$wsGraphs.Range("E4").Formula = "=CORREL(Y,Y_poly)"
$wsGraphs.Range("E5").Formula = "=CORREL(Y,Y_exp)"

$wsGraphs.Range("F2").Formula = '=IF(E2=MAX($E$2:$E$5), "Лучшее", "")'
$wsGraphs.Range("F3").Formula = '=IF(E3=MAX($E$2:$E$5), "Лучшее", "")'
$wsGraphs.Range("F4").Formula = '=IF(E4=MAX($E$2:$E$5), "Лучшее", "")'
$wsGraphs.Range("F5").Formula = '=IF(E5=MAX($E$2:$E$5), "Лучшее", "")'

# --- 4. Active tab / selections --------------------------------------------
$wsData = $wb.Worksheets.Item("Данные")
$wsData.Range("G2:G21").Select()

$wsGraphs.Activate()
$wsGraphs.Range("G3").Select()
